$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 10 was blank; restore it to the "Write program to query database from
# IDE / Noah / 2023-11-28" task (same task as row 6, later deadline).
$ws.Range("A10").Value = "Write program to query database from IDE"
$ws.Range("B10").Value = "Noah"
$ws.Range("C10").Value = 45258

# Drop the "Completed" column (D) that had been added to Table1, shrinking
# the table back down to columns A:C.
$lo = $ws.ListObjects.Item(1)
$lo.ListColumns.Item(4).Delete()

# The column delete above only clears the values/formatting that lived in
# column D; remove the now-empty column outright so the sheet's used range
# collapses back to A2:C12 (dimension, row spans, selection).
$ws.Range("D:D").Delete() | Out-Null

$ws.Range("A2:C12").Select() | Out-Null
